{"js": "// Mill and Refinery address remove cache\n// 1) Mill name -> new name\n// 2) Refinery address -> replacement text\n// 3) \"Given this 3\" + superscript \"rd\" -> \"Given this 10\" + superscript \"th\"\n\nconst body = context.document.body;\n\n// --- 1. Mill / company name (title line) ---\nconst millResults = body.search(\"Busco Sugar Milling Company, Inc.\", { matchCase: true });\nmillResults.load(\"text\");\nawait context.sync();\nif (millResults.items.length > 0) {\n  millResults.items[0].insertText(\n    \"Sweet Crystals Integrated Sugar Mill Corp. - San Fernando\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- 2. Refinery address (italic run) ---\nconst addrResults = body.search(\n  \"4/F Corinthian Plaza Bldg., Paseo de Roxas, Makati City\",\n  { matchCase: true }\n);\naddrResults.load(\"text\");\nawait context.sync();\nif (addrResults.items.length > 0) {\n  addrResults.items[0].insertText(\"aaaaaaaaaaaaaaaaaaaaaaaaaaaaa12123\", \"Replace\");\n  await context.sync();\n}\n\n// --- 3a. \"Given this 3\" -> \"Given this 10\" ---\nconst givenResults = body.search(\"Given this 3\", { matchCase: true });\ngivenResults.load(\"text\");\nawait context.sync();\nif (givenResults.items.length > 0) {\n  givenResults.items[0].insertText(\"Given this 10\", \"Replace\");\n  await context.sync();\n}\n\n// --- 3b. superscript \"rd\" -> \"th\" (the ordinal suffix run right after the date) ---\nconst rdResults = body.search(\"rd\", { matchCase: true });\nrdResults.load(\"text,font\");\nawait context.sync();\nlet ordinalRun = null;\nfor (const item of rdResults.items) {\n  if (item.font.superscript) {\n    ordinalRun = item;\n    break;\n  }\n}\nif (ordinalRun) {\n  ordinalRun.insertText(\"th\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Mill and Refinery address remove cache\n# 1) Mill name -> new name\n# 2) Refinery address -> replacement text\n# 3) \"Given this 3\" + superscript \"rd\" -> \"Given this 10\" + superscript \"th\"\n\n$d = $word.ActiveDocument\n\n# --- 1. Mill / company name (title line) ---\n$find1 = $d.Content.Find\n$find1.Text = \"Busco Sugar Milling Company, Inc.\"\n$find1.Replacement.Text = \"Sweet Crystals Integrated Sugar Mill Corp. - San Fernando\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# --- 2. Refinery address (italic run) ---\n$find2 = $d.Content.Find\n$find2.Text = \"4/F Corinthian Plaza Bldg., Paseo de Roxas, Makati City\"\n$find2.Replacement.Text = \"aaaaaaaaaaaaaaaaaaaaaaaaaaaaa12123\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# --- 3a. \"Given this 3\" -> \"Given this 10\" ---\n$find3 = $d.Content.Find\n$find3.Text = \"Given this 3\"\n$find3.Replacement.Text = \"Given this 10\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n\n# --- 3b. superscript \"rd\" -> \"th\" (the ordinal suffix run right after the date) ---\n# Locate \"Given this 10\" then grab the two characters immediately after it,\n# which is the separate (superscript) run holding the ordinal suffix.\n$rdRange = $d.Content\n$rdRange.Find.Text = \"Given this 10\"\n$rdRange.Find.Execute()\n$rdRange.Collapse(0)\n$rdRange.MoveEnd(1, 2)\n$rdRange.Text = \"th\"\n"}
